$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.322.99'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.94%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.037.67'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +3.96%  '

$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '248.07'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.90%  '

$ws.Range("E6").Value = '  +2.46%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.74'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.54%  '

$ws.Range("E8").Value = '  +0.01%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.400'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +6.77%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0812'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.13%  '

$ws.Range("E11").Value = '  +2.09%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '15.37'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +8.20%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.867'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +4.48%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '2.332.52'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +3.81%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '22.47'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +3.03%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.54'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +5.29%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.031.20'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +3.36%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.220.64'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.99%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '70.92'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.86%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0871'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.18%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.28'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +4.01%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '231.31'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.81%  '

$ws.Range("E23").Value = '  +0.20%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.52'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +3.32%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.37'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +0.93%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.51'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +3.31%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '164.04'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +1.72%  '

$ws.Range("E28").Value = '  -2.92%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '19.90'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +3.01%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.41'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.95%  '

$ws.Range("E31").Value = '  +2.87%  '

$ws.Range("E32").Value = '  +2.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.0670'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +9.20%  '

$ws.Range("E34").Value = '  +2.49%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.53'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +11.73%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.56'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +2.51%  '

$ws.Range("E37").Value = '  -0.06%  '

$ws.Range("E38").Value = '  +2.10%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.43'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.88%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0984'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.79%  '

$ws.Range("E41").Value = '  +1.57%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '17.27'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +9.04%  '

$ws.Range("E43").Value = '  +2.96%  '

$ws.Range("E44").Value = '  +2.93%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '92.93'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.70%  '

$ws.Range("E46").Value = '  +4.30%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.391.08'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.08%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '7.52'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +5.92%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.19'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +20.54%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.87'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +1.21%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '46.80'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.08%  '
